$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:G2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = '257.79'
$ws.Range("E2").Value = '0.00%'
$ws.Range("G2").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D3:G3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = '27.08'
$ws.Range("E3").Value = '-2.13%'
$ws.Range("G3").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D4:G4")
$rng.NumberFormat = "@"
$ws.Range("D4").Value = '4.582'
$ws.Range("E4").Value = '-12.48%'
$ws.Range("G4").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D5:G5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = '0.05886'
$ws.Range("E5").Value = '-0.57%'
$ws.Range("G5").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D6:G6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = '6.643'
$ws.Range("E6").Value = '-0.63%'
$ws.Range("G6").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D7:G7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = '0.8565'
$ws.Range("E7").Value = '-1.38%'
$ws.Range("G7").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D8:G8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = '0.9261'
$ws.Range("E8").Value = '-11.48%'
$ws.Range("G8").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("E9:G9")
$rng.NumberFormat = "@"
$ws.Range("E9").Value = '-0.78%'
$ws.Range("G9").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D10:G10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = '0.03899'
$ws.Range("E10").Value = '7.35%'
$ws.Range("G10").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D11:G11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = '0.07088'
$ws.Range("E11").Value = '-1.36%'
$ws.Range("G11").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D12:G12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = '0.03175'
$ws.Range("E12").Value = '-2.47%'
$ws.Range("G12").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D13:G13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = '0.09180'
$ws.Range("E13").Value = '-0.32%'
$ws.Range("G13").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D14:G14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = '0.001558'
$ws.Range("E14").Value = '0.01%'
$ws.Range("G14").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("B15:G15")
$rng.NumberFormat = "@"
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").Value = '0.0006057'
$ws.Range("E15").Value = '-0.59%'
$ws.Range("G15").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("B16:G16")
$rng.NumberFormat = "@"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005983'
$ws.Range("E16").Value = '2.03%'
$ws.Range("G16").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("B17:G17")
$rng.NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.517'
$ws.Range("E17").Value = '1.01%'
$ws.Range("G17").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("B18:G18")
$rng.NumberFormat = "@"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '3.203'
$ws.Range("E18").Value = '-1.93%'
$ws.Range("G18").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("B19:G19")
$rng.NumberFormat = "@"
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '2.204'
$ws.Range("E19").Value = '-1.00%'
$ws.Range("G19").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D20:G20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = '0.3105'
$ws.Range("E20").Value = '-1.41%'
$ws.Range("G20").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D21:G21")
$rng.NumberFormat = "@"
$ws.Range("D21").Value = '0.1292'
$ws.Range("E21").Value = '0.11%'
$ws.Range("G21").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D22:G22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = '3.889'
$ws.Range("E22").Value = '9.64%'
$ws.Range("G22").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D23:G23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = '0.04197'
$ws.Range("E23").Value = '0.67%'
$ws.Range("G23").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D24:G24")
$rng.NumberFormat = "@"
$ws.Range("D24").Value = '0.001222'
$ws.Range("E24").Value = '0.18%'
$ws.Range("G24").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D25:G25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = '0.004295'
$ws.Range("E25").Value = '-5.32%'
$ws.Range("G25").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("E26:G26")
$rng.NumberFormat = "@"
$ws.Range("E26").Value = '-0.19%'
$ws.Range("G26").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D27:G27")
$rng.NumberFormat = "@"
$ws.Range("D27").Value = '0.0001936'
$ws.Range("E27").Value = '-0.14%'
$ws.Range("G27").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G28")
$rng.NumberFormat = "@"
$ws.Range("G28").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G29")
$rng.NumberFormat = "@"
$ws.Range("G29").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G30")
$rng.NumberFormat = "@"
$ws.Range("G30").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G31")
$rng.NumberFormat = "@"
$ws.Range("G31").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G32")
$rng.NumberFormat = "@"
$ws.Range("G32").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G33")
$rng.NumberFormat = "@"
$ws.Range("G33").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G34")
$rng.NumberFormat = "@"
$ws.Range("G34").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G35")
$rng.NumberFormat = "@"
$ws.Range("G35").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G36")
$rng.NumberFormat = "@"
$ws.Range("G36").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G37")
$rng.NumberFormat = "@"
$ws.Range("G37").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G38")
$rng.NumberFormat = "@"
$ws.Range("G38").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G39")
$rng.NumberFormat = "@"
$ws.Range("G39").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D40:G40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = '0.03851'
$ws.Range("E40").Value = '0.97%'
$ws.Range("G40").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("B41:G41")
$rng.NumberFormat = "@"
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").Value = '0.1103'
$ws.Range("E41").Value = '-0.25%'
$ws.Range("G41").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("B42:G42")
$rng.NumberFormat = "@"
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").Value = '0.003966'
$ws.Range("E42").Value = '-26.90%'
$ws.Range("G42").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D43:G43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = '0.002428'
$ws.Range("E43").Value = '1.93%'
$ws.Range("G43").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D44:G44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = '0.01178'
$ws.Range("E44").Value = '19.16%'
$ws.Range("G44").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D45:G45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = '0.00005454'
$ws.Range("E45").Value = '0.44%'
$ws.Range("G45").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D46:G46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").Value = '-0.08%'
$ws.Range("G46").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D47:G47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = '0.05997'
$ws.Range("E47").Value = '-45.05%'
$ws.Range("G47").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D48:G48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = '0.1263'
$ws.Range("E48").Value = '5,761.88%'
$ws.Range("G48").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D49:G49")
$rng.NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").Value = '-0.08%'
$ws.Range("G49").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("D50:G50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").Value = '-0.08%'
$ws.Range("G50").Value = '20'
$rng.Style = "Normal"

$rng = $ws.Range("G51")
$rng.NumberFormat = "@"
$ws.Range("G51").Value = '20'
$rng.Style = "Normal"
